$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.166777666666666
$ws.Range("H2").Value = 6.500332999999999
$ws.Range("I2").Value = 0.3769884032128669
$ws.Range("J2").Value = 0.376988403212867
$ws.Range("M2").Value = 92.253011
$ws.Range("N2").Value = 276.759033
$ws.Range("O2").Value = 0.2854710184133813
$ws.Range("P2").Value = 0.2854710184133813
$ws.Range("Q2").Value = 199.8917639175543
$ws.Range("R2").Value = 1799.025875257989
$ws.Range("S2").Value = 0.1076192633952116
$ws.Range("T2").Value = 0.1076192633952116

$ws.Range("G3").Value = 2.166777666666666
$ws.Range("H3").Value = 6.500332999999999
$ws.Range("I3").Value = 0.3769884032128669
$ws.Range("J3").Value = 0.376988403212867
$ws.Range("O3").Value = 0.1080355352256351
$ws.Range("P3").Value = 0.1080355352256351
$ws.Range("Q3").Value = 75.64835765835133
$ws.Range("R3").Value = 680.8352189251621
$ws.Range("S3").Value = 0.0407281439149596
$ws.Range("T3").Value = 0.04072814391495961

$ws.Range("G4").Value = 2.166777666666666
$ws.Range("H4").Value = 6.500332999999999
$ws.Range("I4").Value = 0.3769884032128669
$ws.Range("J4").Value = 0.376988403212867
$ws.Range("M4").Value = 42.21774566666667
$ws.Range("N4").Value = 126.653237
$ws.Range("O4").Value = 0.1306401029076487
$ws.Range("P4").Value = 0.1306401029076487
$ws.Range("Q4").Value = 91.47646844754678
$ws.Range("R4").Value = 823.2882160279211
$ws.Range("S4").Value = 0.04924980379071911
$ws.Range("T4").Value = 0.04924980379071912

$ws.Range("G5").Value = 2.166777666666666
$ws.Range("H5").Value = 6.500332999999999
$ws.Range("I5").Value = 0.3769884032128669
$ws.Range("J5").Value = 0.376988403212867
$ws.Range("M5").Value = 13.65158233333333
$ws.Range("N5").Value = 40.954747
$ws.Range("O5").Value = 0.0422439448794879
$ws.Range("P5").Value = 0.0422439448794879
$ws.Range("Q5").Value = 29.57994371452788
$ws.Range("R5").Value = 266.219493430751
$ws.Range("S5").Value = 0.01592547732553051
$ws.Range("T5").Value = 0.01592547732553051

$ws.Range("G6").Value = 2.166777666666666
$ws.Range("H6").Value = 6.500332999999999
$ws.Range("I6").Value = 0.3769884032128669
$ws.Range("J6").Value = 0.376988403212867
$ws.Range("M6").Value = 18.17840666666667
$ws.Range("N6").Value = 54.53522
$ws.Range("O6").Value = 0.0562519120841046
$ws.Range("P6").Value = 0.05625191208410459
$ws.Range("Q6").Value = 39.38856558091778
$ws.Range("R6").Value = 354.49709022826
$ws.Range("S6").Value = 0.02120631851425717
$ws.Range("T6").Value = 0.02120631851425717

$ws.Range("G7").Value = 2.166777666666666
$ws.Range("H7").Value = 6.500332999999999
$ws.Range("I7").Value = 0.3769884032128669
$ws.Range("J7").Value = 0.376988403212867
$ws.Range("M7").Value = 121.9471053333333
$ws.Range("N7").Value = 365.841316
$ws.Range("O7").Value = 0.3773574864897424
$ws.Range("P7").Value = 0.3773574864897424
$ws.Range("Q7").Value = 264.2322643509142
$ws.Range("R7").Value = 2378.090379158228
$ws.Range("S7").Value = 0.142259396272189
$ws.Range("T7").Value = 0.142259396272189

$ws.Range("I8").Value = 0.3757968909097267
$ws.Range("J8").Value = 0.3757968909097268
$ws.Range("M8").Value = 92.253011
$ws.Range("N8").Value = 276.759033
$ws.Range("O8").Value = 0.2854710184133813
$ws.Range("P8").Value = 0.2854710184133813
$ws.Range("Q8").Value = 199.2599845472227
$ws.Range("R8").Value = 1793.339860925004
$ws.Range("S8").Value = 0.1072791211645821
$ws.Range("T8").Value = 0.1072791211645821

$ws.Range("I9").Value = 0.3757968909097267
$ws.Range("J9").Value = 0.3757968909097268
$ws.Range("O9").Value = 0.1080355352256351
$ws.Range("P9").Value = 0.1080355352256351
$ws.Range("S9").Value = 0.04059941824556192
$ws.Range("T9").Value = 0.04059941824556193

$ws.Range("I10").Value = 0.3757968909097267
$ws.Range("J10").Value = 0.3757968909097268
$ws.Range("M10").Value = 42.21774566666667
$ws.Range("N10").Value = 126.653237
$ws.Range("O10").Value = 0.1306401029076487
$ws.Range("P10").Value = 0.1306401029076487
$ws.Range("Q10").Value = 91.18734725263957
$ws.Range("R10").Value = 820.6861252737561
$ws.Range("S10").Value = 0.04909414450082115
$ws.Range("T10").Value = 0.04909414450082115

$ws.Range("I11").Value = 0.3757968909097267
$ws.Range("J11").Value = 0.3757968909097268
$ws.Range("M11").Value = 13.65158233333333
$ws.Range("N11").Value = 40.954747
$ws.Range("O11").Value = 0.0422439448794879
$ws.Range("P11").Value = 0.0422439448794879
$ws.Range("Q11").Value = 29.48645312818178
$ws.Range("R11").Value = 265.378078153636
$ws.Range("S11").Value = 0.01587514314547343
$ws.Range("T11").Value = 0.01587514314547343

$ws.Range("I12").Value = 0.3757968909097267
$ws.Range("J12").Value = 0.3757968909097268
$ws.Range("M12").Value = 18.17840666666667
$ws.Range("N12").Value = 54.53522
$ws.Range("O12").Value = 0.0562519120841046
$ws.Range("P12").Value = 0.05625191208410459
$ws.Range("Q12").Value = 39.26407379259556
$ws.Range("R12").Value = 353.37666413336
$ws.Range("S12").Value = 0.02113929366893379
$ws.Range("T12").Value = 0.0211392936689338

$ws.Range("I13").Value = 0.3757968909097267
$ws.Range("J13").Value = 0.3757968909097268
$ws.Range("M13").Value = 121.9471053333333
$ws.Range("N13").Value = 365.841316
$ws.Range("O13").Value = 0.3773574864897424
$ws.Range("P13").Value = 0.3773574864897424
$ws.Range("Q13").Value = 263.3971299245565
$ws.Range("R13").Value = 2370.574169321008
$ws.Range("S13").Value = 0.1418097701843544
$ws.Range("T13").Value = 0.1418097701843544

$ws.Range("G14").Value = 1.420890666666667
$ws.Range("H14").Value = 4.262672
$ws.Range("I14").Value = 0.2472147058774063
$ws.Range("J14").Value = 0.2472147058774063
$ws.Range("M14").Value = 92.253011
$ws.Range("N14").Value = 276.759033
$ws.Range("O14").Value = 0.2854710184133813
$ws.Range("P14").Value = 0.2854710184133813
$ws.Range("Q14").Value = 131.0814423017974
$ws.Range("R14").Value = 1179.732980716176
$ws.Range("S14").Value = 0.07057263385358771
$ws.Range("T14").Value = 0.0705726338535877

$ws.Range("G15").Value = 1.420890666666667
$ws.Range("H15").Value = 4.262672
$ws.Range("I15").Value = 0.2472147058774063
$ws.Range("J15").Value = 0.2472147058774063
$ws.Range("O15").Value = 0.1080355352256351
$ws.Range("P15").Value = 0.1080355352256351
$ws.Range("Q15").Value = 49.60732566104534
$ws.Range("R15").Value = 446.4659309494081
$ws.Range("S15").Value = 0.02670797306511354
$ws.Range("T15").Value = 0.02670797306511355

$ws.Range("G16").Value = 1.420890666666667
$ws.Range("H16").Value = 4.262672
$ws.Range("I16").Value = 0.2472147058774063
$ws.Range("J16").Value = 0.2472147058774063
$ws.Range("M16").Value = 42.21774566666667
$ws.Range("N16").Value = 126.653237
$ws.Range("O16").Value = 0.1306401029076487
$ws.Range("P16").Value = 0.1306401029076487
$ws.Range("Q16").Value = 59.98680078547379
$ws.Range("R16").Value = 539.8812070692641
$ws.Range("S16").Value = 0.03229615461610848
$ws.Range("T16").Value = 0.03229615461610848

$ws.Range("G17").Value = 1.420890666666667
$ws.Range("H17").Value = 4.262672
$ws.Range("I17").Value = 0.2472147058774063
$ws.Range("J17").Value = 0.2472147058774063
$ws.Range("M17").Value = 13.65158233333333
$ws.Range("N17").Value = 40.954747
$ws.Range("O17").Value = 0.0422439448794879
$ws.Range("P17").Value = 0.0422439448794879
$ws.Range("Q17").Value = 19.39740592266489
$ws.Range("R17").Value = 174.576653303984
$ws.Range("S17").Value = 0.01044332440848397
$ws.Range("T17").Value = 0.01044332440848397

$ws.Range("G18").Value = 1.420890666666667
$ws.Range("H18").Value = 4.262672
$ws.Range("I18").Value = 0.2472147058774063
$ws.Range("J18").Value = 0.2472147058774063
$ws.Range("M18").Value = 18.17840666666667
$ws.Range("N18").Value = 54.53522
$ws.Range("O18").Value = 0.0562519120841046
$ws.Range("P18").Value = 0.05625191208410459
$ws.Range("Q18").Value = 25.82952836753778
$ws.Range("R18").Value = 232.46575530784
$ws.Range("S18").Value = 0.01390629990091364
$ws.Range("T18").Value = 0.01390629990091364

$ws.Range("G19").Value = 1.420890666666667
$ws.Range("H19").Value = 4.262672
$ws.Range("I19").Value = 0.2472147058774063
$ws.Range("J19").Value = 0.2472147058774063
$ws.Range("M19").Value = 121.9471053333333
$ws.Range("N19").Value = 365.841316
$ws.Range("O19").Value = 0.3773574864897424
$ws.Range("P19").Value = 0.3773574864897424
$ws.Range("Q19").Value = 173.2735037951502
$ws.Range("R19").Value = 1559.461534156352
$ws.Range("S19").Value = 0.09328832003319898
$ws.Range("T19").Value = 0.09328832003319899
